$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.760.31'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.037.26'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +3.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.42'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.50'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +2.41%  '

$ws.Range("E7").Value = '  +0.89%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.598'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +3.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.21'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +3.12%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.527.94'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +3.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.63'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.76'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +0.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.058.08'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +3.72%  '

$ws.Range("E17").Value = '  -1.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.54'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -12.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.735.26'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.05'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -0.35%  '

$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("E22").Value = '  +1.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.06'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +0.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.27'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +1.10%  '

$ws.Range("E25").Value = '  -0.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.18'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +1.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.55'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +6.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.174'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +6.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.38'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +3.14%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.32'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +2.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.34'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +2.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.06'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.45'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("E36").Value = '  +5.32%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +8.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.290'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +11.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.15'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +3.42%  '

$ws.Range("E41").Value = '  +3.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +3.28%  '

$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '127.38'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +6.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.76'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +6.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.87'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("E47").Value = '  +5.67%  '

$ws.Range("E48").Value = '  +3.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.039.90'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +1.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.341.28'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0321'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +2.66%  '
